# Auto-generated edits applying the Kujata_Profits.xlsx diff
# (update pricing/profit columns H:N per-row across all 8 job sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 407.36365
$ws.Range("I4").Value = 300.1111
$ws.Range("J4").Value = 890
$ws.Range("K4").Value = 300.1111
$ws.Range("L4").Value = 890
$ws.Range("M4").Value = -186.1111
$ws.Range("N4").Value = -1118
$ws.Range("H12").Value = 180.33333
$ws.Range("I12").Value = 180.33333
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 180.33333
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -10.33332999999999
$ws.Range("N12").ClearContents()
$ws.Range("H28").Value = 144.93333
$ws.Range("I28").Value = 137.72728
$ws.Range("K28").Value = 137.72728
$ws.Range("M28").Value = 347.27272
$ws.Range("H33").Value = 157.45454
$ws.Range("I33").Value = 86
$ws.Range("K33").Value = 86
$ws.Range("M33").Value = 143
$ws.Range("H40").Value = 2698.9167
$ws.Range("I40").Value = 3899.5
$ws.Range("J40").Value = 2098.625
$ws.Range("K40").Value = 3899.5
$ws.Range("L40").Value = 2098.625
$ws.Range("M40").Value = -3724.5
$ws.Range("N40").Value = -2448.625
$ws.Range("H76").Value = 6271.8423
$ws.Range("I76").Value = 6502
$ws.Range("J76").Value = 6228.6875
$ws.Range("K76").Value = 6502
$ws.Range("L76").Value = 6228.6875
$ws.Range("M76").Value = -6187
$ws.Range("N76").Value = -6858.6875
$ws.Range("H79").Value = 6271.8423
$ws.Range("I79").Value = 6502
$ws.Range("J79").Value = 6228.6875
$ws.Range("K79").Value = 6502
$ws.Range("L79").Value = 6228.6875
$ws.Range("M79").Value = -5410
$ws.Range("N79").Value = -8412.6875
$ws.Range("H127").Value = 1063.75
$ws.Range("J127").Value = 2000
$ws.Range("L127").Value = 6000
$ws.Range("N127").Value = -15920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2945.4934
$ws.Range("I32").Value = 2823.5605
$ws.Range("J32").Value = 3839.6667
$ws.Range("K32").Value = 2823.5605
$ws.Range("L32").Value = 3839.6667
$ws.Range("M32").Value = -2536.5605
$ws.Range("N32").Value = -4413.6667
$ws.Range("H92").Value = 1021154
$ws.Range("J92").Value = 1021154
$ws.Range("L92").Value = 1021154
$ws.Range("N92").Value = -1026146
$ws.Range("H102").Value = 23813298
$ws.Range("I102").Value = 27781514
$ws.Range("K102").Value = 27781514
$ws.Range("M102").Value = -27779892
$ws.Range("H114").Value = 20379.3
$ws.Range("J114").Value = 20379.3
$ws.Range("L114").Value = 20379.3
$ws.Range("N114").Value = -29057.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 15625707
$ws.Range("I94").Value = 20834022
$ws.Range("J94").Value = 759.75
$ws.Range("K94").Value = 20834022
$ws.Range("L94").Value = 759.75
$ws.Range("M94").Value = -20833571
$ws.Range("N94").Value = -1661.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 2800
$ws.Range("I17").Value = 2800
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 2800
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -2626
$ws.Range("N17").ClearContents()
$ws.Range("H58").Value = 1579.2727
$ws.Range("I58").Value = 1382.25
$ws.Range("J58").Value = 2104.6667
$ws.Range("K58").Value = 1382.25
$ws.Range("L58").Value = 2104.6667
$ws.Range("M58").Value = -1179.25
$ws.Range("N58").Value = -2510.6667
$ws.Range("H59").Value = 25999.8
$ws.Range("I59").Value = 19999
$ws.Range("J59").Value = 27500
$ws.Range("K59").Value = 19999
$ws.Range("L59").Value = 27500
$ws.Range("M59").Value = -18854
$ws.Range("N59").Value = -29790
$ws.Range("H105").Value = 801.8
$ws.Range("I105").Value = 752.25
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 752.25
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 994.75
$ws.Range("N105").Value = -4494
$ws.Range("H122").Value = 763.7895
$ws.Range("I122").Value = 694.13336
$ws.Range("J122").Value = 1025
$ws.Range("K122").Value = 2082.40008
$ws.Range("L122").Value = 3075
$ws.Range("M122").Value = 367.5999199999997
$ws.Range("N122").Value = -7975
$ws.Range("H134").Value = 1354.1305
$ws.Range("I134").Value = 1221.1904
$ws.Range("K134").Value = 3663.5712
$ws.Range("M134").Value = -1128.5712
$ws.Range("H136").Value = 1579.2727
$ws.Range("I136").Value = 1382.25
$ws.Range("J136").Value = 2104.6667
$ws.Range("K136").Value = 4146.75
$ws.Range("L136").Value = 6314.000100000001
$ws.Range("M136").Value = -1596.75
$ws.Range("N136").Value = -11414.0001
$ws.Range("H141").Value = 27984
$ws.Range("J141").Value = 27984
$ws.Range("L141").Value = 27984
$ws.Range("N141").Value = -38344

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 780541.2
$ws.Range("I4").Value = 100026.555
$ws.Range("J4").Value = 1801313.1
$ws.Range("K4").Value = 300079.665
$ws.Range("L4").Value = 5403939.300000001
$ws.Range("M4").Value = -299967.665
$ws.Range("N4").Value = -5404163.300000001
$ws.Range("H122").Value = 702.1111
$ws.Range("I122").Value = 382.8
$ws.Range("J122").Value = 1101.25
$ws.Range("K122").Value = 3445.2
$ws.Range("L122").Value = 9911.25
$ws.Range("M122").Value = -995.2000000000003
$ws.Range("N122").Value = -14811.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 33999.5
$ws.Range("J103").Value = 33999.5
$ws.Range("L103").Value = 33999.5
$ws.Range("N103").Value = -36343.5
$ws.Range("H113").Value = 1129.45
$ws.Range("I113").Value = 1031
$ws.Range("J113").Value = 1227.9
$ws.Range("K113").Value = 1031
$ws.Range("L113").Value = 1227.9
$ws.Range("M113").Value = 1139
$ws.Range("N113").Value = -5567.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3166.0833
$ws.Range("I40").Value = 2980.6
$ws.Range("K40").Value = 2980.6
$ws.Range("M40").Value = -2844.6
$ws.Range("H122").Value = 50002280
$ws.Range("I122").Value = 62502250
$ws.Range("K122").Value = 187506750
$ws.Range("M122").Value = -187504300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 200
$ws.Range("I9").Value = 200
$ws.Range("K9").Value = 200
$ws.Range("M9").Value = -60
$ws.Range("H70").Value = 15000
$ws.Range("J70").Value = 15000
$ws.Range("L70").Value = 15000
$ws.Range("N70").Value = -15630
$ws.Range("H73").Value = 15000
$ws.Range("J73").Value = 15000
$ws.Range("L73").Value = 15000
$ws.Range("N73").Value = -17184
$ws.Range("H122").Value = 10834650
$ws.Range("I122").Value = 13001240
$ws.Range("K122").Value = 39003720
$ws.Range("M122").Value = -39001270
$ws.Range("H126").Value = 76924550
$ws.Range("I126").Value = 76924550
$ws.Range("K126").Value = 230773650
$ws.Range("M126").Value = -230771180

